$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CVNA")

# Row 4: Inventory
$ws.Range("B4").Value = 1036000000.0
$ws.Range("C4").Value = 968000000.0
$ws.Range("D4").Value = 629000000.0
$ws.Range("E4").Value = 845000000.0
$ws.Range("F4").Value = 763000000.0

# Row 12: Accounts Payable
$ws.Range("B12").Value = 67000000.0
$ws.Range("C12").Value = 96000000.0
$ws.Range("D12").Value = 70000000.0
$ws.Range("E12").Value = 56000000.0
$ws.Range("F12").Value = 64000000.0
